$d = $word.ActiveDocument

$d.Content.Find.Execute("49+48=", $true, $false, $false, $false, $false, $true, 1, $false, "29+47=", 2)
$d.Content.Find.Execute("29+63=", $true, $false, $false, $false, $false, $true, 1, $false, "80-43=", 2)
$d.Content.Find.Execute("61-23=", $true, $false, $false, $false, $false, $true, 1, $false, "62-25=", 2)
$d.Content.Find.Execute("2+49=", $true, $false, $false, $false, $false, $true, 1, $false, "28+57=", 2)
$d.Content.Find.Execute("80-62=", $true, $false, $false, $false, $false, $true, 1, $false, "61-37=", 2)
$d.Content.Find.Execute("24+28=", $true, $false, $false, $false, $false, $true, 1, $false, "12+49=", 2)
$d.Content.Find.Execute("52+19=", $true, $false, $false, $false, $false, $true, 1, $false, "48+7=", 2)
$d.Content.Find.Execute("41-38=", $true, $false, $false, $false, $false, $true, 1, $false, "28+29=", 2)
$d.Content.Find.Execute("81-36=", $true, $false, $false, $false, $false, $true, 1, $false, "91-48=", 2)
$d.Content.Find.Execute("4+39=", $true, $false, $false, $false, $false, $true, 1, $false, "86-58=", 2)
$d.Content.Find.Execute("74+17=", $true, $false, $false, $false, $false, $true, 1, $false, "93-56=", 2)
$d.Content.Find.Execute("81-4=", $true, $false, $false, $false, $false, $true, 1, $false, "78+17=", 2)
$d.Content.Find.Execute("57+5=", $true, $false, $false, $false, $false, $true, 1, $false, "19+72=", 2)
$d.Content.Find.Execute("85+8=", $true, $false, $false, $false, $false, $true, 1, $false, "58-19=", 2)
$d.Content.Find.Execute("22-7=", $true, $false, $false, $false, $false, $true, 1, $false, "58+5=", 2)
$d.Content.Find.Execute("63-17=", $true, $false, $false, $false, $false, $true, 1, $false, "54-16=", 2)
$d.Content.Find.Execute("86+9=", $true, $false, $false, $false, $false, $true, 1, $false, "56-8=", 2)
$d.Content.Find.Execute("9+77=", $true, $false, $false, $false, $false, $true, 1, $false, "68+3=", 2)
$d.Content.Find.Execute("38+5=", $true, $false, $false, $false, $false, $true, 1, $false, "52-29=", 2)
$d.Content.Find.Execute("40-33=", $true, $false, $false, $false, $false, $true, 1, $false, "9+38=", 2)
$d.Content.Find.Execute("95-67=", $true, $false, $false, $false, $false, $true, 1, $false, "75-37=", 2)
$d.Content.Find.Execute("90-9=", $true, $false, $false, $false, $false, $true, 1, $false, "53-25=", 2)
$d.Content.Find.Execute("93-38=", $true, $false, $false, $false, $false, $true, 1, $false, "77-68=", 2)
$d.Content.Find.Execute("53-16=", $true, $false, $false, $false, $false, $true, 1, $false, "59+16=", 2)
$d.Content.Find.Execute("18+59=", $true, $false, $false, $false, $false, $true, 1, $false, "24+47=", 2)
$d.Content.Find.Execute("19+42=", $true, $false, $false, $false, $false, $true, 1, $false, "81-42=", 2)
$d.Content.Find.Execute("5+69=", $true, $false, $false, $false, $false, $true, 1, $false, "60-13=", 2)
$d.Content.Find.Execute("39+24=", $true, $false, $false, $false, $false, $true, 1, $false, "30-3=", 2)
$d.Content.Find.Execute("19+73=", $true, $false, $false, $false, $false, $true, 1, $false, "42-38=", 2)
$d.Content.Find.Execute("90-18=", $true, $false, $false, $false, $false, $true, 1, $false, "47+27=", 2)
$d.Content.Find.Execute("77-49=", $true, $false, $false, $false, $false, $true, 1, $false, "52-36=", 2)
$d.Content.Find.Execute("51-7=", $true, $false, $false, $false, $false, $true, 1, $false, "74-29=", 2)
$d.Content.Find.Execute("85+9=", $true, $false, $false, $false, $false, $true, 1, $false, "9+48=", 2)
$d.Content.Find.Execute("8+87=", $true, $false, $false, $false, $false, $true, 1, $false, "23+68=", 2)
$d.Content.Find.Execute("30-18=", $true, $false, $false, $false, $false, $true, 1, $false, "7+75=", 2)
$d.Content.Find.Execute("62-35=", $true, $false, $false, $false, $false, $true, 1, $false, "93-57=", 2)
$d.Content.Find.Execute("16+47=", $true, $false, $false, $false, $false, $true, 1, $false, "42-16=", 2)
$d.Content.Find.Execute("77+14=", $true, $false, $false, $false, $false, $true, 1, $false, "70-16=", 2)
$d.Content.Find.Execute("72-48=", $true, $false, $false, $false, $false, $true, 1, $false, "96-59=", 2)
$d.Content.Find.Execute("14+78=", $true, $false, $false, $false, $false, $true, 1, $false, "14+47=", 2)
$d.Content.Find.Execute("49+12=", $true, $false, $false, $false, $false, $true, 1, $false, "84-77=", 2)
$d.Content.Find.Execute("67+8=", $true, $false, $false, $false, $false, $true, 1, $false, "48+6=", 2)
$d.Content.Find.Execute("81-79=", $true, $false, $false, $false, $false, $true, 1, $false, "36+19=", 2)
$d.Content.Find.Execute("78-29=", $true, $false, $false, $false, $false, $true, 1, $false, "36+26=", 2)
$d.Content.Find.Execute("20-7=", $true, $false, $false, $false, $false, $true, 1, $false, "37+34=", 2)
$d.Content.Find.Execute("43+9=", $true, $false, $false, $false, $false, $true, 1, $false, "86-67=", 2)
$d.Content.Find.Execute("48+48=", $true, $false, $false, $false, $false, $true, 1, $false, "89+3=", 2)
$d.Content.Find.Execute("44-35=", $true, $false, $false, $false, $false, $true, 1, $false, "39+48=", 2)
$d.Content.Find.Execute("84-58=", $true, $false, $false, $false, $false, $true, 1, $false, "92-46=", 2)
$d.Content.Find.Execute("62-47=", $true, $false, $false, $false, $false, $true, 1, $false, "73-58=", 2)
$d.Content.Find.Execute("92-43=", $true, $false, $false, $false, $false, $true, 1, $false, "36+46=", 2)
$d.Content.Find.Execute("82-33=", $true, $false, $false, $false, $false, $true, 1, $false, "72-38=", 2)
$d.Content.Find.Execute("73-69=", $true, $false, $false, $false, $false, $true, 1, $false, "63-25=", 2)
$d.Content.Find.Execute("82-9=", $true, $false, $false, $false, $false, $true, 1, $false, "61-58=", 2)
$d.Content.Find.Execute("39+13=", $true, $false, $false, $false, $false, $true, 1, $false, "7+47=", 2)
$d.Content.Find.Execute("92-13=", $true, $false, $false, $false, $false, $true, 1, $false, "66-29=", 2)
$d.Content.Find.Execute("91-17=", $true, $false, $false, $false, $false, $true, 1, $false, "28+19=", 2)
$d.Content.Find.Execute("39+49=", $true, $false, $false, $false, $false, $true, 1, $false, "26-9=", 2)
$d.Content.Find.Execute("18+58=", $true, $false, $false, $false, $false, $true, 1, $false, "70-49=", 2)
$d.Content.Find.Execute("14+8=", $true, $false, $false, $false, $false, $true, 1, $false, "90-52=", 2)
$d.Content.Find.Execute("84-15=", $true, $false, $false, $false, $false, $true, 1, $false, "22-9=", 2)
$d.Content.Find.Execute("28+17=", $true, $false, $false, $false, $false, $true, 1, $false, "94-5=", 2)
$d.Content.Find.Execute("24+57=", $true, $false, $false, $false, $false, $true, 1, $false, "18+14=", 2)
$d.Content.Find.Execute("94-67=", $true, $false, $false, $false, $false, $true, 1, $false, "4+7=", 2)
$d.Content.Find.Execute("6+18=", $true, $false, $false, $false, $false, $true, 1, $false, "56+36=", 2)
$d.Content.Find.Execute("9+17=", $true, $false, $false, $false, $false, $true, 1, $false, "72-65=", 2)
$d.Content.Find.Execute("61-46=", $true, $false, $false, $false, $false, $true, 1, $false, "39+48=", 2)
$d.Content.Find.Execute("4+9=", $true, $false, $false, $false, $false, $true, 1, $false, "7+85=", 2)
$d.Content.Find.Execute("73-28=", $true, $false, $false, $false, $false, $true, 1, $false, "36+29=", 2)
$d.Content.Find.Execute("17+18=", $true, $false, $false, $false, $false, $true, 1, $false, "86-79=", 2)
$d.Content.Find.Execute("55+6=", $true, $false, $false, $false, $false, $true, 1, $false, "70-7=", 2)
$d.Content.Find.Execute("14+79=", $true, $false, $false, $false, $false, $true, 1, $false, "18+43=", 2)
$d.Content.Find.Execute("6+37=", $true, $false, $false, $false, $false, $true, 1, $false, "51-29=", 2)
$d.Content.Find.Execute("81-23=", $true, $false, $false, $false, $false, $true, 1, $false, "61-54=", 2)
$d.Content.Find.Execute("38+35=", $true, $false, $false, $false, $false, $true, 1, $false, "65+6=", 2)
$d.Content.Find.Execute("70-69=", $true, $false, $false, $false, $false, $true, 1, $false, "57+37=", 2)
$d.Content.Find.Execute("72-43=", $true, $false, $false, $false, $false, $true, 1, $false, "94-56=", 2)
$d.Content.Find.Execute("59+27=", $true, $false, $false, $false, $false, $true, 1, $false, "61-12=", 2)
$d.Content.Find.Execute("59+32=", $true, $false, $false, $false, $false, $true, 1, $false, "85-7=", 2)
$d.Content.Find.Execute("51-45=", $true, $false, $false, $false, $false, $true, 1, $false, "91-3=", 2)
$d.Content.Find.Execute("54+18=", $true, $false, $false, $false, $false, $true, 1, $false, "15-6=", 2)
$d.Content.Find.Execute("39+42=", $true, $false, $false, $false, $false, $true, 1, $false, "68+6=", 2)
$d.Content.Find.Execute("72-3=", $true, $false, $false, $false, $false, $true, 1, $false, "71-66=", 2)
$d.Content.Find.Execute("83-79=", $true, $false, $false, $false, $false, $true, 1, $false, "39+5=", 2)
$d.Content.Find.Execute("64-18=", $true, $false, $false, $false, $false, $true, 1, $false, "47-8=", 2)
$d.Content.Find.Execute("80-33=", $true, $false, $false, $false, $false, $true, 1, $false, "84-57=", 2)
$d.Content.Find.Execute("17+14=", $true, $false, $false, $false, $false, $true, 1, $false, "38+45=", 2)
$d.Content.Find.Execute("15+46=", $true, $false, $false, $false, $false, $true, 1, $false, "63+29=", 2)
$d.Content.Find.Execute("46+38=", $true, $false, $false, $false, $false, $true, 1, $false, "60-26=", 2)
$d.Content.Find.Execute("70-22=", $true, $false, $false, $false, $false, $true, 1, $false, "19+34=", 2)
$d.Content.Find.Execute("61-39=", $true, $false, $false, $false, $false, $true, 1, $false, "17+76=", 2)
$d.Content.Find.Execute("57+35=", $true, $false, $false, $false, $false, $true, 1, $false, "35+48=", 2)
$d.Content.Find.Execute("49+4=", $true, $false, $false, $false, $false, $true, 1, $false, "58+26=", 2)
$d.Content.Find.Execute("53+19=", $true, $false, $false, $false, $false, $true, 1, $false, "93-45=", 2)
$d.Content.Find.Execute("14+38=", $true, $false, $false, $false, $false, $true, 1, $false, "59+15=", 2)
$d.Content.Find.Execute("9+4=", $true, $false, $false, $false, $false, $true, 1, $false, "96-37=", 2)
$d.Content.Find.Execute("75-59=", $true, $false, $false, $false, $false, $true, 1, $false, "10-5=", 2)
$d.Content.Find.Execute("66+9=", $true, $false, $false, $false, $false, $true, 1, $false, "70-63=", 2)
$d.Content.Find.Execute("42-19=", $true, $false, $false, $false, $false, $true, 1, $false, "54-25=", 2)
